# "cambio de formato plantilla excel de productos"
# Turns the single header row into a two-tier header:
#   - a new row 1 marking each column "obligatorio" (required) or "opcional" (optional)
#   - row 2 left blank (spacer)
#   - the original header labels pushed down to row 3, keeping their bold/border style
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing header row (row 1) down to row 3, inserting two fresh rows
# above it. Excel's row-insert carries the original cell styles/content along,
# so "nombre/codigo/..." + their s="1" (and G's s="2") formatting lands on row 3
# automatically.
$ws.Rows("1:2").Insert()

# Column A now holds real content ("nombre") for the first time - size it to fit,
# matching the bestFit columns already on the sheet.
$ws.Columns("A").AutoFit()

# New row 1: mark which columns are required vs optional for the import template.
$ws.Range("A1").Value2 = "obligatorio"
$ws.Range("B1").Value2 = "obligatorio"
$ws.Range("C1:G1").Value2 = "opcional"

# G3 ("stock actual") previously carried a redundant applyFill style distinct from
# the rest of the header row; normalize it back to the same look as its neighbors
# (fill was already "none" either way, so this is purely a formatting cleanup).
$ws.Range("G3").Interior.Pattern = -4142

# Leave the cursor where the author ended up after editing the template.
$ws.Range("F9").Select() | Out-Null
